$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.87
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.55
$ws.Range("J3").Value = 2.42
$ws.Range("K3").Value = 2.02
$ws.Range("N3").Value = 6
$ws.Range("U3").Value = 1.95
$ws.Range("W3").Value = 5.8
$ws.Range("AA3").Value = 16.5
$ws.Range("AC3").Value = 6
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 16
$ws.Range("AF3").Value = 90
$ws.Range("AH3").Value = 10.5
$ws.Range("AJ3").Value = 15
$ws.Range("AL3").Value = 50
$ws.Range("AM3").Value = 55
$ws.Range("AN3").Value = 3.65
$ws.Range("AQ3").Value = 35
$ws.Range("AR3").Value = 70
$ws.Range("AU3").Value = 7.2
$ws.Range("AX3").Value = 27
$ws.Range("AY3").Value = 32
$ws.Range("BA3").Value = 200
$ws.Range("BB3").Value = 450

# Row 4 updates
$ws.Range("H4").Value = 2.6
$ws.Range("U4").Value = 1.93
$ws.Range("V4").Value = 1.78
$ws.Range("AD4").Value = 5.1
$ws.Range("AG4").Value = 800
$ws.Range("AH4").Value = 6.7
$ws.Range("AQ4").Value = 100
$ws.Range("AU4").Value = 6.7

# Row 8 updates
$ws.Range("G8").Value = 4.1
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 1.8
$ws.Range("J8").Value = 4.5
$ws.Range("L8").Value = 2.4
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 8.5
$ws.Range("W8").Value = 12
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 15
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 7
$ws.Range("AE8").Value = 17
$ws.Range("AI8").Value = 8.5
$ws.Range("AL8").Value = 15
$ws.Range("AN8").Value = 6
$ws.Range("AO8").Value = 23
$ws.Range("AQ8").Value = 81
$ws.Range("AU8").Value = 8.5
$ws.Range("AW8").Value = 3.75
$ws.Range("AX8").Value = 9.5
$ws.Range("AZ8").Value = 29
